$wb = $excel.ActiveWorkbook

# Sheet index 1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(12, 8).Value = 1069.6  # H12: 907.8333 -> 1069.6
$ws.Cells.Item(12, 9).Value = 832.5  # I12: 727.7143 -> 832.5
$ws.Cells.Item(12, 10).Value = 1425.25  # J12: 1160 -> 1425.25
$ws.Cells.Item(12, 11).Value = 832.5  # K12: 727.7143 -> 832.5
$ws.Cells.Item(12, 12).Value = 1425.25  # L12: 1160 -> 1425.25
$ws.Cells.Item(12, 13).Value = -662.5  # M12: -557.7143 -> -662.5
$ws.Cells.Item(12, 14).Value = -1765.25  # N12: -1500 -> -1765.25
$ws.Cells.Item(100, 8).Value = 1306.2858  # H100: 1312.2142 -> 1306.2858
$ws.Cells.Item(100, 9).Value = 806.9091  # I100: 781.25 -> 806.9091
$ws.Cells.Item(100, 10).Value = 1629.4117  # J100: 1710.4375 -> 1629.4117
$ws.Cells.Item(100, 11).Value = 806.9091  # K100: 781.25 -> 806.9091
$ws.Cells.Item(100, 12).Value = 1629.4117  # L100: 1710.4375 -> 1629.4117
$ws.Cells.Item(100, 13).Value = -265.9091  # M100: -240.25 -> -265.9091
$ws.Cells.Item(100, 14).Value = -2711.4117  # N100: -2792.4375 -> -2711.4117
$ws.Cells.Item(116, 8).Value = 16347  # H116: 10673.5 -> 16347
$ws.Cells.Item(116, 10).Value = 20166.666  # J116: 11500 -> 20166.666
$ws.Cells.Item(116, 12).Value = 20166.666  # L116: 11500 -> 20166.666
$ws.Cells.Item(116, 14).Value = -27050.666  # N116: -18384 -> -27050.666

# Sheet index 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 19578630  # H2: 20002052 -> 19578630
$ws.Cells.Item(2, 9).Value = 22983224  # I2: 23811544 -> 22983224
$ws.Cells.Item(2, 11).Value = 22983224  # K2: 23811544 -> 22983224
$ws.Cells.Item(2, 13).Value = -22983111  # M2: -23811431 -> -22983111
$ws.Cells.Item(5, 8).Value = 310.5  # H5: 366.6 -> 310.5
$ws.Cells.Item(5, 9).Value = 190.75  # I5: 244.33333 -> 190.75
$ws.Cells.Item(5, 11).Value = 190.75  # K5: 244.33333 -> 190.75
$ws.Cells.Item(5, 13).Value = -78.75  # M5: -132.33333 -> -78.75
$ws.Cells.Item(40, 8).Value = 34749.5  # H40: 20000 -> 34749.5
$ws.Cells.Item(40, 9).Value = 49499  # I40: 0 -> 49499
$ws.Cells.Item(40, 11).Value = 49499  # K40: 0 -> 49499
$ws.Cells.Item(40, 13).Value = -49323  # M40: None -> -49323
$ws.Cells.Item(45, 8).Value = 11354.934  # H45: 12840.462 -> 11354.934
$ws.Cells.Item(45, 9).Value = 14678.4  # I45: 16120.556 -> 14678.4
$ws.Cells.Item(45, 10).Value = 4708  # J45: 5460.25 -> 4708
$ws.Cells.Item(45, 11).Value = 14678.4  # K45: 16120.556 -> 14678.4
$ws.Cells.Item(45, 12).Value = 4708  # L45: 5460.25 -> 4708
$ws.Cells.Item(45, 13).Value = -14301.4  # M45: -15743.556 -> -14301.4
$ws.Cells.Item(45, 14).Value = -5462  # N45: -6214.25 -> -5462
$ws.Cells.Item(49, 8).Value = 0  # H49: 5000 -> 0
$ws.Cells.Item(49, 10).Value = 0  # J49: 5000 -> 0
$ws.Cells.Item(49, 12).Value = 0  # L49: 5000 -> 0
$ws.Cells.Item(49, 14).ClearContents()  # N49: -5520 -> (removed)
$ws.Cells.Item(97, 8).Value = 33373934  # H97: 38507620 -> 33373934
$ws.Cells.Item(97, 9).Value = 41673980  # I97: 50007780 -> 41673980
$ws.Cells.Item(97, 11).Value = 41673980  # K97: 50007780 -> 41673980
$ws.Cells.Item(97, 13).Value = -41673484  # M97: -50007284 -> -41673484
$ws.Cells.Item(110, 8).Value = 1597.2128  # H110: 1657.4419 -> 1597.2128
$ws.Cells.Item(110, 9).Value = 1687.3954  # I110: 1763.0513 -> 1687.3954
$ws.Cells.Item(110, 11).Value = 1687.3954  # K110: 1763.0513 -> 1687.3954
$ws.Cells.Item(110, 13).Value = 357.6045999999999  # M110: 281.9486999999999 -> 357.6045999999999
$ws.Cells.Item(116, 8).Value = 19578630  # H116: 20002052 -> 19578630
$ws.Cells.Item(116, 9).Value = 22983224  # I116: 23811544 -> 22983224
$ws.Cells.Item(116, 11).Value = 22983224  # K116: 23811544 -> 22983224
$ws.Cells.Item(116, 13).Value = -22980930  # M116: -23809250 -> -22980930

# Sheet index 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 19578630  # H3: 20002052 -> 19578630
$ws.Cells.Item(3, 9).Value = 22983224  # I3: 23811544 -> 22983224
$ws.Cells.Item(3, 11).Value = 22983224  # K3: 23811544 -> 22983224
$ws.Cells.Item(3, 13).Value = -22983110  # M3: -23811430 -> -22983110
$ws.Cells.Item(4, 8).Value = 310.5  # H4: 366.6 -> 310.5
$ws.Cells.Item(4, 9).Value = 190.75  # I4: 244.33333 -> 190.75
$ws.Cells.Item(4, 11).Value = 190.75  # K4: 244.33333 -> 190.75
$ws.Cells.Item(4, 13).Value = -75.75  # M4: -129.33333 -> -75.75
$ws.Cells.Item(105, 8).Value = 4367.5938  # H105: 4229.4165 -> 4367.5938
$ws.Cells.Item(105, 9).Value = 3120.16  # I105: 3123.2144 -> 3120.16
$ws.Cells.Item(105, 10).Value = 8822.714  # J105: 8101.125 -> 8822.714
$ws.Cells.Item(105, 11).Value = 3120.16  # K105: 3123.2144 -> 3120.16
$ws.Cells.Item(105, 12).Value = 8822.714  # L105: 8101.125 -> 8822.714
$ws.Cells.Item(105, 13).Value = -1373.16  # M105: -1376.2144 -> -1373.16
$ws.Cells.Item(105, 14).Value = -12316.714  # N105: -11595.125 -> -12316.714
$ws.Cells.Item(134, 8).Value = 7829.5454  # H134: 9653.471 -> 7829.5454
$ws.Cells.Item(134, 9).Value = 7862.55  # I134: 9940.666999999999 -> 7862.55
$ws.Cells.Item(134, 11).Value = 23587.65  # K134: 29822.001 -> 23587.65
$ws.Cells.Item(134, 13).Value = -21052.65  # M134: -27287.001 -> -21052.65

# Sheet index 4
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(62, 8).Value = 200006180  # H62: 250005860 -> 200006180
$ws.Cells.Item(62, 9).Value = 250005730  # I62: 333338460 -> 250005730
$ws.Cells.Item(62, 11).Value = 250005730  # K62: 333338460 -> 250005730
$ws.Cells.Item(62, 13).Value = -250005106  # M62: -333337836 -> -250005106
$ws.Cells.Item(65, 8).Value = 200006180  # H65: 250005860 -> 200006180
$ws.Cells.Item(65, 9).Value = 250005730  # I65: 333338460 -> 250005730
$ws.Cells.Item(65, 11).Value = 1250028650  # K65: 1666692300 -> 1250028650
$ws.Cells.Item(65, 13).Value = -1250025530  # M65: -1666689180 -> -1250025530

# Sheet index 5
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 8).Value = 71.708336  # H2: 76.90909000000001 -> 71.708336
$ws.Cells.Item(2, 9).Value = 27.692308  # I2: 30.1 -> 27.692308
$ws.Cells.Item(2, 10).Value = 123.72727  # J2: 115.916664 -> 123.72727
$ws.Cells.Item(2, 11).Value = 166.153848  # K2: 180.6 -> 166.153848
$ws.Cells.Item(2, 12).Value = 742.3636200000001  # L2: 695.499984 -> 742.3636200000001
$ws.Cells.Item(2, 13).Value = -53.15384800000001  # M2: -67.60000000000002 -> -53.15384800000001
$ws.Cells.Item(2, 14).Value = -968.3636200000001  # N2: -921.499984 -> -968.3636200000001
$ws.Cells.Item(4, 8).Value = 118439.766  # H4: 118439.94 -> 118439.766
$ws.Cells.Item(4, 9).Value = 83456.836  # I4: 83457.086 -> 83456.836
$ws.Cells.Item(4, 11).Value = 250370.508  # K4: 250371.258 -> 250370.508
$ws.Cells.Item(4, 13).Value = -250258.508  # M4: -250259.258 -> -250258.508
$ws.Cells.Item(12, 8).Value = 83.333336  # H12: 99.69231000000001 -> 83.333336
$ws.Cells.Item(12, 9).Value = 68.666664  # I12: 7 -> 68.666664
$ws.Cells.Item(12, 10).Value = 87  # J12: 107.416664 -> 87
$ws.Cells.Item(12, 11).Value = 205.999992  # K12: 21 -> 205.999992
$ws.Cells.Item(12, 12).Value = 261  # L12: 322.249992 -> 261
$ws.Cells.Item(12, 13).Value = -32.99999199999999  # M12: 152 -> -32.99999199999999
$ws.Cells.Item(12, 14).Value = -607  # N12: -668.249992 -> -607
$ws.Cells.Item(16, 8).Value = 200  # H16: 350 -> 200
$ws.Cells.Item(16, 10).Value = 0  # J16: 500 -> 0
$ws.Cells.Item(16, 12).Value = 0  # L16: 1500 -> 0
$ws.Cells.Item(16, 14).ClearContents()  # N16: -1846 -> (removed)
$ws.Cells.Item(19, 8).Value = 300  # H19: 400 -> 300
$ws.Cells.Item(19, 10).Value = 300  # J19: 400 -> 300
$ws.Cells.Item(19, 12).Value = 900  # L19: 1200 -> 900
$ws.Cells.Item(19, 14).Value = -1248  # N19: -1548 -> -1248
$ws.Cells.Item(104, 8).Value = 19791.438  # H104: 19799.8 -> 19791.438
$ws.Cells.Item(104, 10).Value = 19974.285  # J104: 19998 -> 19974.285
$ws.Cells.Item(104, 12).Value = 59922.855  # L104: 59994 -> 59922.855
$ws.Cells.Item(104, 14).Value = -65164.855  # N104: -65236 -> -65164.855
$ws.Cells.Item(118, 8).Value = 625  # H118: 600 -> 625
$ws.Cells.Item(118, 9).Value = 625  # I118: 600 -> 625
$ws.Cells.Item(118, 11).Value = 1875  # K118: 1800 -> 1875
$ws.Cells.Item(118, 13).Value = -632  # M118: -557 -> -632
$ws.Cells.Item(140, 8).Value = 1105.1364  # H140: 1139.238 -> 1105.1364
$ws.Cells.Item(140, 9).Value = 1105.1364  # I140: 1139.238 -> 1105.1364
$ws.Cells.Item(140, 11).Value = 3315.4092  # K140: 3417.714 -> 3315.4092
$ws.Cells.Item(140, 13).Value = 1864.5908  # M140: 1762.286 -> 1864.5908

# Sheet index 6
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(41, 8).Value = 21384.572  # H41: 8365.5 -> 21384.572
$ws.Cells.Item(41, 9).Value = 23615.334  # I41: 8438.6 -> 23615.334
$ws.Cells.Item(41, 11).Value = 23615.334  # K41: 8438.6 -> 23615.334
$ws.Cells.Item(41, 13).Value = -23260.334  # M41: -8083.6 -> -23260.334
$ws.Cells.Item(97, 8).Value = 1352  # H97: 1269.5714 -> 1352
$ws.Cells.Item(97, 9).Value = 1016.25  # I97: 968 -> 1016.25
$ws.Cells.Item(97, 11).Value = 1016.25  # K97: 968 -> 1016.25
$ws.Cells.Item(97, 13).Value = -520.25  # M97: -472 -> -520.25
$ws.Cells.Item(122, 8).Value = 3510.7778  # H122: 3713.125 -> 3510.7778
$ws.Cells.Item(122, 9).Value = 3760.1  # I122: 4227.125 -> 3760.1
$ws.Cells.Item(122, 11).Value = 11280.3  # K122: 12681.375 -> 11280.3
$ws.Cells.Item(122, 13).Value = -8830.299999999999  # M122: -10231.375 -> -8830.299999999999

# Sheet index 7
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(16, 8).Value = 9618050  # H16: 12503196 -> 9618050
$ws.Cells.Item(16, 9).Value = 12502886  # I16: 13891995 -> 12502886
$ws.Cells.Item(16, 10).Value = 1932.3334  # J16: 4002 -> 1932.3334
$ws.Cells.Item(16, 11).Value = 12502886  # K16: 13891995 -> 12502886
$ws.Cells.Item(16, 12).Value = 1932.3334  # L16: 4002 -> 1932.3334
$ws.Cells.Item(16, 13).Value = -12502716  # M16: -13891825 -> -12502716
$ws.Cells.Item(16, 14).Value = -2272.3334  # N16: -4342 -> -2272.3334
$ws.Cells.Item(50, 8).Value = 40000  # H50: 0 -> 40000
$ws.Cells.Item(50, 10).Value = 40000  # J50: 0 -> 40000
$ws.Cells.Item(50, 12).Value = 40000  # L50: 0 -> 40000
$ws.Cells.Item(50, 14).Value = -41274  # N50: None -> -41274
$ws.Cells.Item(74, 8).Value = 43798.75  # H74: 47548.5 -> 43798.75
$ws.Cells.Item(74, 9).Value = 37598.5  # I74: 40197 -> 37598.5
$ws.Cells.Item(74, 11).Value = 37598.5  # K74: 40197 -> 37598.5
$ws.Cells.Item(74, 13).Value = -36600.5  # M74: -39199 -> -36600.5
$ws.Cells.Item(77, 8).Value = 43798.75  # H77: 47548.5 -> 43798.75
$ws.Cells.Item(77, 9).Value = 37598.5  # I77: 40197 -> 37598.5
$ws.Cells.Item(77, 11).Value = 112795.5  # K77: 120591 -> 112795.5
$ws.Cells.Item(77, 13).Value = -107803.5  # M77: -115599 -> -107803.5
$ws.Cells.Item(82, 8).Value = 454.07  # H82: 458.4898 -> 454.07
$ws.Cells.Item(82, 9).Value = 406.42554  # I82: 408.37634 -> 406.42554
$ws.Cells.Item(82, 10).Value = 1200.5  # J82: 1390.6 -> 1200.5
$ws.Cells.Item(82, 11).Value = 406.42554  # K82: 408.37634 -> 406.42554
$ws.Cells.Item(82, 12).Value = 1200.5  # L82: 1390.6 -> 1200.5
$ws.Cells.Item(82, 13).Value = -45.42554000000001  # M82: -47.37634000000003 -> -45.42554000000001
$ws.Cells.Item(82, 14).Value = -1922.5  # N82: -2112.6 -> -1922.5
$ws.Cells.Item(85, 8).Value = 454.07  # H85: 458.4898 -> 454.07
$ws.Cells.Item(85, 9).Value = 406.42554  # I85: 408.37634 -> 406.42554
$ws.Cells.Item(85, 10).Value = 1200.5  # J85: 1390.6 -> 1200.5
$ws.Cells.Item(85, 11).Value = 406.42554  # K85: 408.37634 -> 406.42554
$ws.Cells.Item(85, 12).Value = 1200.5  # L85: 1390.6 -> 1200.5
$ws.Cells.Item(85, 13).Value = 841.57446  # M85: 839.62366 -> 841.57446
$ws.Cells.Item(85, 14).Value = -3696.5  # N85: -3886.6 -> -3696.5
$ws.Cells.Item(93, 8).Value = 7408300  # H93: 7693257.5 -> 7408300
$ws.Cells.Item(93, 9).Value = 8696497  # I93: 9091819 -> 8696497
$ws.Cells.Item(93, 11).Value = 8696497  # K93: 9091819 -> 8696497
$ws.Cells.Item(93, 13).Value = -8695249  # M93: -9090571 -> -8695249
$ws.Cells.Item(122, 8).Value = 7640  # H122: 7436.029 -> 7640
$ws.Cells.Item(122, 9).Value = 7740.9565  # I122: 7287.385 -> 7740.9565
$ws.Cells.Item(122, 10).Value = 7407.8  # J122: 7865.4443 -> 7407.8
$ws.Cells.Item(122, 11).Value = 23222.8695  # K122: 21862.155 -> 23222.8695
$ws.Cells.Item(122, 12).Value = 22223.4  # L122: 23596.3329 -> 22223.4
$ws.Cells.Item(122, 13).Value = -20772.8695  # M122: -19412.155 -> -20772.8695
$ws.Cells.Item(122, 14).Value = -27123.4  # N122: -28496.3329 -> -27123.4
$ws.Cells.Item(136, 8).Value = 6435504  # H136: 6437004 -> 6435504
$ws.Cells.Item(136, 10).Value = 11544.75  # J136: 15044.583 -> 11544.75
$ws.Cells.Item(136, 12).Value = 34634.25  # L136: 45133.749 -> 34634.25
$ws.Cells.Item(136, 14).Value = -39734.25  # N136: -50233.749 -> -39734.25

# Sheet index 8
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(34, 8).Value = 49498.5  # H34: 49749 -> 49498.5
$ws.Cells.Item(34, 9).Value = 49498  # I34: 49999 -> 49498
$ws.Cells.Item(34, 11).Value = 49498  # K34: 49999 -> 49498
$ws.Cells.Item(34, 13).Value = -49295  # M34: -49796 -> -49295
$ws.Cells.Item(37, 8).Value = 0  # H37: 49999 -> 0
$ws.Cells.Item(37, 9).Value = 0  # I37: 49999 -> 0
$ws.Cells.Item(37, 11).Value = 0  # K37: 49999 -> 0
$ws.Cells.Item(37, 13).ClearContents()  # M37: -49796 -> (removed)
$ws.Cells.Item(43, 8).Value = 34005  # H43: 34905 -> 34005
$ws.Cells.Item(43, 9).Value = 34005  # I43: 34905 -> 34005
$ws.Cells.Item(43, 11).Value = 34005  # K43: 34905 -> 34005
$ws.Cells.Item(43, 13).Value = -33856  # M43: -34756 -> -33856
$ws.Cells.Item(107, 8).Value = 23810028  # H107: 33333852 -> 23810028
$ws.Cells.Item(107, 9).Value = 587.3333  # I107: 638 -> 587.3333
$ws.Cells.Item(107, 10).Value = 166666670  # J107: 83333670 -> 166666670
$ws.Cells.Item(107, 11).Value = 1761.9999  # K107: 1914 -> 1761.9999
$ws.Cells.Item(107, 12).Value = 500000010  # L107: 250001010 -> 500000010
$ws.Cells.Item(107, 13).Value = 158.0001  # M107: 6 -> 158.0001
$ws.Cells.Item(107, 14).Value = -500003850  # N107: -250004850 -> -500003850
$ws.Cells.Item(113, 8).Value = 481.55554  # H113: 495.80768 -> 481.55554
$ws.Cells.Item(113, 9).Value = 309.8889  # I113: 321.58823 -> 309.8889
$ws.Cells.Item(113, 11).Value = 929.6667  # K113: 964.76469 -> 929.6667
$ws.Cells.Item(113, 13).Value = 1240.3333  # M113: 1205.23531 -> 1240.3333
$ws.Cells.Item(126, 8).Value = 8957.52  # H126: 9348.48 -> 8957.52
$ws.Cells.Item(126, 9).Value = 7997.952  # I126: 8186.75 -> 7997.952
$ws.Cells.Item(126, 10).Value = 13995.25  # J126: 13995.4 -> 13995.25
$ws.Cells.Item(126, 11).Value = 23993.856  # K126: 24560.25 -> 23993.856
$ws.Cells.Item(126, 12).Value = 41985.75  # L126: 41986.2 -> 41985.75
$ws.Cells.Item(126, 13).Value = -21523.856  # M126: -22090.25 -> -21523.856
$ws.Cells.Item(126, 14).Value = -46925.75  # N126: -46926.2 -> -46925.75
$ws.Cells.Item(136, 8).Value = 1543.5625  # H136: 1589.2413 -> 1543.5625
$ws.Cells.Item(136, 9).Value = 1339.4231  # I136: 1412.9166 -> 1339.4231
$ws.Cells.Item(136, 10).Value = 2428.1667  # J136: 2435.6 -> 2428.1667
$ws.Cells.Item(136, 11).Value = 4018.2693  # K136: 4238.7498 -> 4018.2693
$ws.Cells.Item(136, 12).Value = 7284.500100000001  # L136: 7306.799999999999 -> 7284.500100000001
$ws.Cells.Item(136, 13).Value = -1468.2693  # M136: -1688.7498 -> -1468.2693
$ws.Cells.Item(136, 14).Value = -12384.5001  # N136: -12406.8 -> -12384.5001
